# Edit: add new "tn.5250" terminal-emulation command category.
#  - Inserts a new column before Z on the hidden '#system' sheet and
#    populates it with the new "tn.5250" profile commands.
#  - Inserts a new "ocr(image,saveVar)" command row inside the existing
#    "image" column (K) and renames the "colorbit" command's first
#    parameter from "source" to "image".
#  - Inserts a new "tn.5250" row inside the existing "target" column (A).
#  - Adds two new top-level image commands to the shared-strings table
#    indirectly via the two new sheet rows above.
#  - Updates every definedName so the named ranges keep pointing at the
#    right (now shifted) columns, and adds the new "tn.5250" defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new blank column before column Z (shifts Z:AE -> AA:AF)
# ---------------------------------------------------------------------
$ws.Columns("Z:Z").Insert()

# ---------------------------------------------------------------------
# 2) Insert a new row inside column K only (the "image" command list),
#    right after "crop(...)" / before "resize(...)".
#    NOTE: Range.Insert() on a single-column range shifts the whole
#    sheet row in this engine (unlike real Excel's single-column cell
#    insert), so instead we push K6:K7 down to K7:K8 with plain value
#    writes and then drop the new command into the vacated K6.
# ---------------------------------------------------------------------
$ws.Range("K8").Value = $ws.Range("K7").Value
$ws.Range("K7").Value = $ws.Range("K6").Value
$ws.Range("K6").Value = "ocr(image,saveVar)"

# Rename colorbit's first parameter from "source" to "image"
$ws.Range("K2").Value = "colorbit(image,bit,saveTo)"

# ---------------------------------------------------------------------
# 3) Insert a new row inside column A only (the "target" list), right
#    before "web", and set its value to the new category name. Same
#    column-only shift restriction applies here, so push the tail of
#    the list down with value writes instead of Range.Insert().
# ---------------------------------------------------------------------
$ws.Range("A32").Value = $ws.Range("A31").Value
$ws.Range("A31").Value = $ws.Range("A30").Value
$ws.Range("A30").Value = $ws.Range("A29").Value
$ws.Range("A29").Value = $ws.Range("A28").Value
$ws.Range("A28").Value = $ws.Range("A27").Value
$ws.Range("A27").Value = $ws.Range("A26").Value
$ws.Range("A26").Value = "tn.5250"

# ---------------------------------------------------------------------
# 4) Populate the new column Z with the "tn.5250" header + commands.
# ---------------------------------------------------------------------
$ws.Range("Z1").Value = "tn.5250"
$ws.Range("Z2").Value = "close(profile)"
$ws.Range("Z3").Value = "open(profile)"
$ws.Range("Z4").Value = "saveText(profile,var)"
$ws.Range("Z5").Value = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value = "updateScreenFields(profile)"

# ---------------------------------------------------------------------
# 5) Fix up the defined names: ranges that live at/after column Z need
#    to point one column to the right, "target" and "image" now cover
#    one extra row, and "tn.5250" is a brand-new name.
# ---------------------------------------------------------------------
$wb.Names.Item("image").RefersTo    = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo   = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo      = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo= "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo       = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo      = "='#system'!`$AF`$2:`$AF`$27"

$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
